$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 170, shifting existing rows 170-201 down to 171-202.
$ws.Rows(170).Insert()

# Populate the newly inserted row 170 with the new data record.
$ws.Cells.Item(170, 1).Value = 10
$ws.Cells.Item(170, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(170, 3).Value = "La Araucanía"
$ws.Cells.Item(170, 4).Value = 45015
$ws.Cells.Item(170, 5).Value = 9
$ws.Cells.Item(170, 6).Value = 100112031
$ws.Cells.Item(170, 7).Value = "Poroto verde"
$ws.Cells.Item(170, 8).Value = "Brío"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 300
$ws.Cells.Item(170, 11).Value = 1000
$ws.Cells.Item(170, 12).Value = 1000
$ws.Cells.Item(170, 13).Value = 1000
$ws.Cells.Item(170, 14).Value = "$/kilo"
$ws.Cells.Item(170, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(170, 16).Value = 1000
$ws.Cells.Item(170, 17).Value = 1
$ws.Cells.Item(170, 18).Value = "Hortaliza"
